# Auto-generated edit script: apply 2025-06-25 violent-crime data increment
# Updates 2025 (column L) year-to-date totals, plus a few prior-column
# corrections, across the Citywide Totals, By Neighborhood, and per-
# neighborhood sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3093
$ws.Range("L3").Value = 3155
$ws.Range("F4").Value = 1929
$ws.Range("J4").Value = 1867
$ws.Range("L4").Value = 813
$ws.Range("L5").Value = 178
$ws.Range("L6").Value = 2804
$ws.Range("F7").Value = 24122
$ws.Range("J7").Value = 29342
$ws.Range("L7").Value = 10043

$ws = $wb.Sheets.Item("Grant Park")
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 9

$ws = $wb.Sheets.Item("Logan Square")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 113

$ws = $wb.Sheets.Item("Austin")
$ws.Range("L2").Value = 187
$ws.Range("L3").Value = 211
$ws.Range("L7").Value = 640

$ws = $wb.Sheets.Item("South Chicago")
$ws.Range("L3").Value = 95
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 238

$ws = $wb.Sheets.Item("Garfield Park")
$ws.Range("L3").Value = 143
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 464

$ws = $wb.Sheets.Item("West Pullman")
$ws.Range("L2").Value = 54
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 136

$ws = $wb.Sheets.Item("New City")
$ws.Range("L2").Value = 70
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 189

$ws = $wb.Sheets.Item("By Neighborhood")
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 336
$ws.Range("L8").Value = 640
$ws.Range("L10").Value = 66
$ws.Range("L11").Value = 167
$ws.Range("L15").Value = 74
$ws.Range("L19").Value = 283
$ws.Range("L20").Value = 254
$ws.Range("L27").Value = 97
$ws.Range("L29").Value = 546
$ws.Range("L33").Value = 464
$ws.Range("L36").Value = 138
$ws.Range("L38").Value = 9
$ws.Range("L42").Value = 321
$ws.Range("L45").Value = 18
$ws.Range("L50").Value = 52
$ws.Range("L51").Value = 123
$ws.Range("L52").Value = 199
$ws.Range("L53").Value = 113
$ws.Range("L54").Value = 206
$ws.Range("L58").Value = 6
$ws.Range("L60").Value = 62
$ws.Range("F63").Value = 213
$ws.Range("J63").Value = 218
$ws.Range("L63").Value = 34
$ws.Range("L64").Value = 67
$ws.Range("L65").Value = 189
$ws.Range("L67").Value = 366
$ws.Range("L71").Value = 31
$ws.Range("L76").Value = 135
$ws.Range("L77").Value = 60
$ws.Range("L78").Value = 123
$ws.Range("L80").Value = 32
$ws.Range("L81").Value = 11
$ws.Range("L83").Value = 238
$ws.Range("L84").Value = 101
$ws.Range("L85").Value = 507
$ws.Range("L86").Value = 73
$ws.Range("L89").Value = 136
$ws.Range("L91").Value = 144
$ws.Range("L94").Value = 119
$ws.Range("L95").Value = 136
$ws.Range("L96").Value = 99
$ws.Range("L98").Value = 62
$ws.Range("F101").Value = 24122
$ws.Range("J101").Value = 29342
$ws.Range("L101").Value = 10043

$ws = $wb.Sheets.Item("North Lawndale")
$ws.Range("L2").Value = 110
$ws.Range("L3").Value = 136
$ws.Range("L5").Value = 10
$ws.Range("L7").Value = 366

$ws = $wb.Sheets.Item("South Deering")
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 101

$ws = $wb.Sheets.Item("Loop")
$ws.Range("L2").Value = 43
$ws.Range("L6").Value = 102
$ws.Range("L7").Value = 206

$ws = $wb.Sheets.Item("Englewood")
$ws.Range("L2").Value = 166
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 546

$ws = $wb.Sheets.Item("Chatham")
$ws.Range("L3").Value = 86
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 283

$ws = $wb.Sheets.Item("River North")
$ws.Range("L3").Value = 26
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 135

$ws = $wb.Sheets.Item("Ashburn")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 78

$ws = $wb.Sheets.Item("Humboldt Park")
$ws.Range("L2").Value = 94
$ws.Range("L3").Value = 101
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 321

$ws = $wb.Sheets.Item("Avondale")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 66

$ws = $wb.Sheets.Item("Rogers Park")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 123

$ws = $wb.Sheets.Item("West Ridge")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 99

$ws = $wb.Sheets.Item("Washington Park")
$ws.Range("L3").Value = 57
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 144

$ws = $wb.Sheets.Item("Near South Side")
$ws.Range("L2").Value = 23
$ws.Range("L3").Value = 16
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 67

$ws = $wb.Sheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 81
$ws.Range("L7").Value = 254

$ws = $wb.Sheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 138

$ws = $wb.Sheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 103
$ws.Range("L7").Value = 336

$ws = $wb.Sheets.Item("West Loop")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 119

$ws = $wb.Sheets.Item("Brighton Park")
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 74

$ws = $wb.Sheets.Item("Wicker Park")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 62

$ws = $wb.Sheets.Item("Lincoln Square")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 52

$ws = $wb.Sheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 167

$ws = $wb.Sheets.Item("Uptown")
$ws.Range("L2").Value = 43
$ws.Range("L7").Value = 136

$ws = $wb.Sheets.Item("Edgewater")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 97

$ws = $wb.Sheets.Item("Streeterville")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 73

$ws = $wb.Sheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 123

$ws = $wb.Sheets.Item("Morgan Park")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 62

$ws = $wb.Sheets.Item("South Shore")
$ws.Range("L2").Value = 149
$ws.Range("L3").Value = 207
$ws.Range("L6").Value = 103
$ws.Range("L7").Value = 507

$ws = $wb.Sheets.Item("Oakland")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 31

$ws = $wb.Sheets.Item("Riverdale")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 60

$ws = $wb.Sheets.Item("Jackson Park")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 18

$ws = $wb.Sheets.Item("Rush & Division")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 32

$ws = $wb.Sheets.Item("Little Village")
$ws.Range("L2").Value = 70
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 199

$ws = $wb.Sheets.Item("Sauganash,Forest Glen")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 11

$ws = $wb.Sheets.Item("Millenium Park")
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 6

